$d = $word.ActiveDocument

function New-RunXmlPackage([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1) After "14:00: Bereich in Umgebung sperren" insert three new log entries
#    ("14:45: Klo", "15:15 Fortsetzung mit Sperren der Umgebung",
#    "16:00 Fehler......"), pushing the _GoBack bookmark down into its own
#    trailing paragraph (matches how Word relocates a bookmark when a
#    paragraph it sits at the end of is split by newly typed paragraphs).
# ---------------------------------------------------------------------------
$ellipsis = [string][char]0x2026 + [char]0x2026
$findRng = $d.Content
$replaced = $findRng.Find.Execute(
    "14:00: Bereich in Umgebung sperren",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "14:00: Bereich in Umgebung sperren^p14:45: Klo^p15:15 Fortsetzung mit Sperren der Umgebung^p16:00 Fehler$ellipsis^p",
    2)
Write-Host "Step1 replace:" $replaced

# Split the "15:15 Fortsetzung mit Sperren der Umgebung" run into two runs
# ("15:15" and " Fortsetzung mit Sperren der Umgebung"), as in the source.
$splitRng = $d.Content
$foundSplit = $splitRng.Find.Execute("15:15 Fortsetzung mit Sperren der Umgebung")
Write-Host "Step1 find split target:" $foundSplit
$s = $splitRng.Start
$e = $splitRng.End
$freshSplitRng = $d.Range($s, $e)
$splitXml = New-RunXmlPackage('<w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/></w:rPr><w:t>15:15</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Fortsetzung mit Sperren der Umgebung</w:t></w:r></w:p>')
$freshSplitRng.InsertXML($splitXml)
Write-Host "Step1 split done"

# ---------------------------------------------------------------------------
# 2) Move the <w:lastRenderedPageBreak/> that sits before "Programmablaufplan
#    erstellen" to just before "Feld selektieren" instead.
# ---------------------------------------------------------------------------
$feldRng = $d.Content
$foundFeld = $feldRng.Find.Execute("Feld selektieren")
Write-Host "Step2 find Feld selektieren:" $foundFeld
$s = $feldRng.Start
$e = $feldRng.End
$freshFeldRng = $d.Range($s, $e)
$feldXml = New-RunXmlPackage('<w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>Feld selektieren</w:t></w:r></w:p>')
$freshFeldRng.InsertXML($feldXml)
Write-Host "Step2 insert done"

$papRng = $d.Content
$foundPap = $papRng.Find.Execute("Programmablaufplan erstellen")
Write-Host "Step2 find Programmablaufplan erstellen:" $foundPap
$s = $papRng.Start
$e = $papRng.End
$freshPapRng = $d.Range($s, $e)
$papXml = New-RunXmlPackage('<w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:highlight w:val="green"/></w:rPr><w:t>Programmablaufplan erstellen</w:t></w:r></w:p>')
$freshPapRng.InsertXML($papXml)
Write-Host "Step2 remove done"

# ---------------------------------------------------------------------------
# 3) Move the <w:lastRenderedPageBreak/> that sits before "4.6
#    Qualitaetssicherung" to just before "3.5 Qualitaetsanforderungen "
#    instead.
# ---------------------------------------------------------------------------
$qaRng = $d.Content
$foundQa = $qaRng.Find.Execute("3.5 Qualitätsanforderungen ")
Write-Host "Step3 find 3.5 Qualitaetsanforderungen:" $foundQa
$s = $qaRng.Start
$e = $qaRng.End
$freshQaRng = $d.Range($s, $e)
$qaXml = New-RunXmlPackage('<w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">3.5 Qualitätsanforderungen </w:t></w:r></w:p>')
$freshQaRng.InsertXML($qaXml)
Write-Host "Step3 insert done"

$qsRng = $d.Content
$foundQs = $qsRng.Find.Execute("4.6 Qualitätssicherung")
Write-Host "Step3 find 4.6 Qualitaetssicherung:" $foundQs
$s = $qsRng.Start
$e = $qsRng.End
$freshQsRng = $d.Range($s, $e)
$qsXml = New-RunXmlPackage('<w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/></w:rPr><w:t>4.6 Qualitätssicherung</w:t></w:r></w:p>')
$freshQsRng.InsertXML($qsXml)
Write-Host "Step3 remove done"

# ---------------------------------------------------------------------------
# 4) New stackoverflow-link paragraph, inserted into the (until now empty)
#    paragraph right before the "loop-for-each-over-an-array" link.
# ---------------------------------------------------------------------------
$jsonRng = $d.Content
$foundJson = $jsonRng.Find.Execute("json-decode.php")
Write-Host "Step4 find json-decode.php:" $foundJson
$jsonParaRange = $d.Range($jsonRng.Start, $jsonRng.End)
$jsonPara = $jsonParaRange.Paragraphs(1)
$emptyPara = $jsonPara.Next()
$s = $emptyPara.Range.Start
$e = $emptyPara.Range.End - 1
$freshEmptyRng = $d.Range($s, $e)
$findIfXml = New-RunXmlPackage('<w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/></w:rPr><w:t>https://stackoverflow.com/questions/42526032/how-to-find-if-element-with-specific-id-exists-or-not</w:t></w:r></w:p>')
$freshEmptyRng.InsertXML($findIfXml)
Write-Host "Step4 insert done"

# ---------------------------------------------------------------------------
# 5) Remove the stray <w:lastRenderedPageBreak/> before the "increment
#    letters" stackoverflow link (not moved elsewhere).
# ---------------------------------------------------------------------------
$incRng = $d.Content
$foundInc = $incRng.Find.Execute("https://stackoverflow.com/questions/12504042/what-is-a-method-that-can-be-used-to-increment-letters")
Write-Host "Step5 find increment-letters:" $foundInc
$s = $incRng.Start
$e = $incRng.End
$freshIncRng = $d.Range($s, $e)
$incXml = New-RunXmlPackage('<w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/></w:rPr><w:t>https://stackoverflow.com/questions/12504042/what-is-a-method-that-can-be-used-to-increment-letters</w:t></w:r></w:p>')
$freshIncRng.InsertXML($incXml)
Write-Host "Step5 remove done"
